$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing table (A1:C9) already carries a wrap-text / center / top
# cell style at the column level, so plain value assignment on the new
# rows picks up matching formatting automatically - append three more
# rows of test-plan data below the current last row (row 9).

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Test Deliverables"
$ws.Range("C10").Value = "Test Cases Documents,bug Report,Final Summary"

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Testing Tools"
$ws.Range("C11").Value = "Browser dev tools,Postman,Excel"

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Environments Requeriments"
$ws.Range("C12").Value = "Crome,MongoDB,Spring Boot"

# Match the wrapped-text row heights Excel would have auto-fit for this
# column width (20.89 chars) given each row's text length.
$ws.Rows("10").RowHeight = 43.2
$ws.Rows("11").RowHeight = 28.8
$ws.Rows("12").RowHeight = 28.8

$ws.Range("A13").Select() | Out-Null
